$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data refresh rotates the "Fruta / Granada" rows: the same
# underlying records are now reported against different price dates. This
# performs the row-content permutation described by the diff (dates,
# variety, quality, volume, prices, unit, origin, $/kg, kg/unit) while
# leaving rows 7, 8 and 10 untouched.

$rows = @{
    2  = @{ D = 44266; K = "Wonderfull";       L = "Segunda"; M = 120; N = 4800;   O = 4800;   P = 4800;   Q = "`$/bandeja 4 kilos";       R = "Provincia del Elquí";   S = 1200; T = 4 }
    3  = @{ D = 44266; K = "Wonderfull";       L = "Tercera"; M = 80;  N = 4000;   O = 4000;   P = 4000;   Q = "`$/bandeja 4 kilos";       R = "Provincia del Elquí";   S = 1000; T = 4 }
    4  = @{ D = 44285; K = "Wonderfull";       L = "Primera"; M = 8;   N = 280000; O = 300000; P = 290000; Q = "`$/bins (400 kilos)";      R = "Provincia del Elquí";   S = 725;  T = 400 }
    5  = @{ D = 44320; K = "Wonderfull";       L = "Primera"; M = 12;  N = 250000; O = 260000; P = 255000; Q = "`$/bins (400 kilos)";      R = "Provincia de Limarí";   S = 638;  T = 400 }
    6  = @{ D = 44721; K = "Wonderfull";       L = "Primera"; M = 7;   N = 300000; O = 300000; P = 300000; Q = "`$/bins (400 kilos)";      R = "Región Metropolitana";  S = 750;  T = 400 }
    9  = @{ D = 44280; K = "Sin especificar";  L = "Primera"; M = 15;  N = 360000; O = 360000; P = 360000; Q = "`$/bins (450 kilos)";      R = "Provincia del Elquí";   S = 800;  T = 450 }
    11 = @{ D = 44662; K = "Sin especificar";  L = "Primera"; M = 45;  N = 18000;  O = 18000;  P = 18000;  Q = "`$/caja 18 kilos granel";  R = "Provincia de Limarí";   S = 1000; T = 18 }
    12 = @{ D = 44662; K = "Sin especificar";  L = "Segunda"; M = 60;  N = 16000;  O = 16000;  P = 16000;  Q = "`$/caja 18 kilos granel";  R = "Provincia de Limarí";   S = 889;  T = 18 }
    13 = @{ D = 44307; K = "Sin especificar";  L = "Primera"; M = 150; N = 16000;  O = 18000;  P = 17000;  Q = "`$/caja 15 kilos granel";  R = "Región de O'Higgins";   S = 1133; T = 15 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value2 = $vals.D    # D: Fecha
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
}
